$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric/percent-looking text values are written as literal text
# (matching the source inline-string cells) rather than being auto-coerced
# into numbers/percentages by Excel's input parser.
$textCells = @('D2', 'E2', 'D3', 'E3', 'D4', 'E4', 'D5', 'E5', 'D6', 'E6', 'D7', 'E7', 'D8', 'E8', 'D9', 'E9', 'D10', 'E10', 'D11', 'E11', 'D12', 'E12', 'D13', 'E13', 'D14', 'E14', 'D15', 'E15', 'D16', 'E16', 'D17', 'E17', 'D18', 'E18', 'E19', 'E20', 'E21', 'D22', 'E22', 'D23', 'E23', 'D24', 'E24', 'D25', 'E25', 'D26', 'E26', 'D27', 'E27', 'D40', 'E40', 'D41', 'E41', 'D42', 'E42', 'D43', 'E43', 'D44', 'E44', 'D45', 'E45', 'D46', 'E46')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '275.17'
$ws.Range('E2').Value = '-2.55%'
$ws.Range('D3').Value = '27.13'
$ws.Range('E3').Value = '1.24%'
$ws.Range('D4').Value = '4.757'
$ws.Range('E4').Value = '-3.70%'
$ws.Range('D5').Value = '0.06298'
$ws.Range('E5').Value = '-1.78%'
$ws.Range('D6').Value = '6.935'
$ws.Range('E6').Value = '-0.89%'
$ws.Range('B7').Value = 'FTXToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D7').Value = '1.333'
$ws.Range('E7').Value = '26.82%'
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D8').Value = '0.8766'
$ws.Range('E8').Value = '-1.06%'
$ws.Range('B9').Value = 'WazirX'
$ws.Range('C9').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D9').Value = '0.1508'
$ws.Range('E9').Value = '0.78%'
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').Value = '0.05066'
$ws.Range('E10').Value = '-1.85%'
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').Value = '0.07502'
$ws.Range('E11').Value = '0.56%'
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').Value = '0.02906'
$ws.Range('E12').Value = '-6.28%'
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D13').Value = '0.09002'
$ws.Range('E13').Value = '-0.51%'
$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D14').Value = '0.001562'
$ws.Range('E14').Value = '-0.54%'
$ws.Range('B15').Value = 'One'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D15').Value = '0.0006373'
$ws.Range('E15').Value = '0.93%'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').Value = '0.006041'
$ws.Range('E16').Value = '-0.05%'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').Value = '3.448'
$ws.Range('E17').Value = '-1.75%'
$ws.Range('B18').Value = 'GateToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D18').Value = '3.302'
$ws.Range('E18').Value = '-1.43%'
$ws.Range('E19').Value = '-1.18%'
$ws.Range('E20').Value = '0.77%'
$ws.Range('E21').Value = '2.61%'
$ws.Range('D22').Value = '3.922'
$ws.Range('E22').Value = '-0.75%'
$ws.Range('D23').Value = '0.04402'
$ws.Range('E23').Value = '1.37%'
$ws.Range('D24').Value = '0.001172'
$ws.Range('E24').Value = '-0.23%'
$ws.Range('D25').Value = '0.003833'
$ws.Range('E25').Value = '3.88%'
$ws.Range('D26').Value = '0.0001201'
$ws.Range('E26').Value = '0.28%'
$ws.Range('D27').Value = '0.0001937'
$ws.Range('E27').Value = '14.58%'
$ws.Range('D40').Value = '0.04100'
$ws.Range('E40').Value = '-0.06%'
$ws.Range('D41').Value = '0.006789'
$ws.Range('E41').Value = '2.14%'
$ws.Range('D42').Value = '0.1171'
$ws.Range('E42').Value = '-0.77%'
$ws.Range('D43').Value = '0.002191'
$ws.Range('E43').Value = '-6.94%'
$ws.Range('D44').Value = '0.01153'
$ws.Range('E44').Value = '-11.90%'
$ws.Range('D45').Value = '0.00005166'
$ws.Range('E45').Value = '-1.61%'
$ws.Range('D46').Value = '0.02301'
$ws.Range('E46').Value = '2.46%'
